$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Links")
Write-Output $ws.Name
